# feat: add 2022-Q1 data
#
# 1) Insert a new leading row into "总计" summarising the new quarter and
#    push the existing rows down by one.
# 2) Insert a new worksheet named "2022-Q1" right before "总计" and fill it
#    with the per-fund holding detail for the new quarter.
#
# NOTE: Worksheets.Add() shifts what "current sheet" means for handles that
# were captured earlier, so every "总计" edit below happens *before* the new
# sheet is added, and every "2022-Q1" edit happens via handles fetched
# *after* Add() returns.

$wb = $excel.ActiveWorkbook

function Set-TextCell($rng, $text) {
    $rng.NumberFormat = "@"
    $rng.Value2 = $text
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. "总计" sheet: push existing rows down one and insert the 2022-Q1 summary
#    at the top. Write bottom-up so source rows aren't clobbered before
#    they're copied down.
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")

# old row 4 (2021-Q2) -> new row 5
Set-TextCell $zongji.Range("B5") "2021-Q2"
$zongji.Range("C5").Value2 = 7
$zongji.Range("D5").Value2 = 0.06

# old row 3 (2021-Q3) -> new row 4
Set-TextCell $zongji.Range("B4") "2021-Q3"
$zongji.Range("C4").Value2 = 12
$zongji.Range("D4").Value2 = 0.87

# old row 2 (2021-Q4) -> new row 3
Set-TextCell $zongji.Range("B3") "2021-Q4"
$zongji.Range("C3").Value2 = 9
$zongji.Range("D3").Value2 = 0.32

# new row 2 (2022-Q1)
Set-TextCell $zongji.Range("B2") "2022-Q1"
$zongji.Range("C2").Value2 = 5
$zongji.Range("D2").Value2 = 0.1

# Extend the bordered/bold "index" style already used on A2:A4 down to the
# new A5, then renumber the whole 0..3 index column.
$zongji.Range("A4").Copy()
$zongji.Range("A5").PasteSpecial(-4122)
$zongji.Range("A2").Value2 = 0
$zongji.Range("A3").Value2 = 1
$zongji.Range("A4").Value2 = 2
$zongji.Range("A5").Value2 = 3

# ---------------------------------------------------------------------------
# 2. New "2022-Q1" sheet, positioned immediately before "总计"
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$added = $wb.Worksheets.Add($zongji)
$added.Name = "2022-Q1"

$q1 = $wb.Worksheets.Item("2022-Q1")

# Header row (row 1) -- write text, then copy the bold/bordered format that
# the sibling quarter sheets already use on B1:H1.
$q1.Range("B1").Value2 = "基金代码"
$q1.Range("C1").Value2 = "基金名称"
$q1.Range("D1").Value2 = "基金规模"
$q1.Range("E1").Value2 = "股票总仓位"
$q1.Range("F1").Value2 = "仓位占比"
$q1.Range("G1").Value2 = "持有市值(亿元)"
$q1.Range("H1").Value2 = "仓位排名"

$prior = $wb.Worksheets.Item("2021-Q4")
$prior.Range("B1:H1").Copy()
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Range("B1:H1").PasteSpecial(-4122)

# Row 2
$q1.Range("A2").Value2 = 0
Set-TextCell $q1.Range("B2") "002415"
Set-TextCell $q1.Range("C2") "融通通盈灵活配置混合"
Set-TextCell $q1.Range("D2") "0.89"
Set-TextCell $q1.Range("E2") "69.03"
Set-TextCell $q1.Range("F2") "6.66"
Set-TextCell $q1.Range("G2") "0.0593"
$q1.Range("H2").Value2 = 2

# Row 3
$q1.Range("A3").Value2 = 1
Set-TextCell $q1.Range("B3") "006700"
Set-TextCell $q1.Range("C3") "红土创新稳健混合A"
Set-TextCell $q1.Range("D3") "0.74"
Set-TextCell $q1.Range("E3") "27.06"
Set-TextCell $q1.Range("F3") "3.82"
Set-TextCell $q1.Range("G3") "0.0283"
$q1.Range("H3").Value2 = 3

# Row 4
$q1.Range("A4").Value2 = 2
Set-TextCell $q1.Range("B4") "006701"
Set-TextCell $q1.Range("C4") "红土创新稳健混合C"
Set-TextCell $q1.Range("D4") "0.35"
Set-TextCell $q1.Range("E4") "27.06"
Set-TextCell $q1.Range("F4") "3.82"
Set-TextCell $q1.Range("G4") "0.0134"
$q1.Range("H4").Value2 = 3

# Row 5
$q1.Range("A5").Value2 = 3
Set-TextCell $q1.Range("B5") "006231"
Set-TextCell $q1.Range("C5") "国融融君灵活配置混合A"
Set-TextCell $q1.Range("D5") "0.08"
Set-TextCell $q1.Range("E5") "57.39"
Set-TextCell $q1.Range("F5") "2.09"
Set-TextCell $q1.Range("G5") "0.0017"
$q1.Range("H5").Value2 = 7

# Row 6
$q1.Range("A6").Value2 = 4
Set-TextCell $q1.Range("B6") "006232"
Set-TextCell $q1.Range("C6") "国融融君灵活配置混合C"
Set-TextCell $q1.Range("D6") "0.00"
Set-TextCell $q1.Range("E6") "57.39"
Set-TextCell $q1.Range("F6") "2.09"
$q1.Range("G6").Value2 = 0
$q1.Range("H6").Value2 = 7

# Column A (index) cells share the bordered/bold style used in row 1; copy
# it down from the sibling sheet's equivalent column.
$prior = $wb.Worksheets.Item("2021-Q4")
$prior.Range("A2:A6").Copy()
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Range("A2:A6").PasteSpecial(-4122)
$q1.Range("A2").Value2 = 0
$q1.Range("A3").Value2 = 1
$q1.Range("A4").Value2 = 2
$q1.Range("A5").Value2 = 3
$q1.Range("A6").Value2 = 4
